$d = $word.ActiveDocument

# The document currently ends with a paragraph that only contains an
# inline picture. Add a brand new paragraph after it with the day's
# closing note, written as two separate runs (the second one carries
# a leading space, so it needs xml:space="preserve").
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$target = $newPara.Range

$xmlFrag = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t>So I started implementing a new AXI controller. I was not using the register feature of AXI devices properly. When trying to edit the files a bit, Vivado kept crashing.</w:t></w:r>
            <w:r><w:t xml:space="preserve"> Gonna close up for today.</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$target.InsertXML($xmlFrag)

# InsertXML drops the new paragraph's markup in place of $target but
# leaves the paragraph mark that used to belong to $newPara dangling as
# an extra empty paragraph right after it. Delete that leftover
# paragraph mark so the new paragraph becomes the final paragraph again,
# exactly as in the target document.
$textPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$trailingMark = $d.Range($textPara.Range.End - 1, $textPara.Range.End)
$trailingMark.Delete()
